# Supplier template: rename telp -> wa column, update phone/WA numbers,
# update address text, and clear two trailing cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: supplier_telp -> supplier_wa ---
$ws.Range("C1").Value = "supplier_wa"
$ws.Range("C1").NumberFormat = "#,##0"

# --- Row 2: update WA number (C2), apply right-aligned number format ---
$ws.Range("C2").Value = 6287866711494
$ws.Range("C2").NumberFormat = "#,##0"
$ws.Range("C2").HorizontalAlignment = -4152
$ws.Range("C2").Font.Name = "Calibri"
$ws.Range("C2").Font.Size = 11

# --- Row 3: update supplier name / address text, clear phone + alamat ---
$ws.Range("B3").Value = "random - pasar kasin"
$ws.Range("C3").NumberFormat = "#,##0"
$ws.Range("C3").HorizontalAlignment = -4152
$ws.Range("C3").Font.Name = "Calibri"
$ws.Range("C3").Font.Size = 11
$ws.Range("C3").Value = "'"
$ws.Range("D3").Value = "'"

# --- widen column C to fit the longer WA header / values ---
$ws.Columns("C").ColumnWidth = 48.14785714285715
